$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update account / IBAN / amount text that changed in this edit.
$ws.Range("A3").Value = "97-1818181-4/QQQ ŞUBESİ"
$ws.Range("A4").Value = "TRQQ 0020 3000 0180 YYYY 0000 OO"

$ws.Range("A9").Value = "97-1818181-4/QQQ ŞUBESİ"
$ws.Range("A10").Value = "TRQQ 0020 3000 0180 YYYY 0000 OO"

$ws.Range("A11").Value = "97-1818181-5/QQQ ŞUBESİ"
$ws.Range("A12").Value = "TRQQ 0020 3000 0180 YYYY 0000 OO"

$ws.Range("E9").Value = "50 USD"
$ws.Range("E11").Value = "50 EUR"

# The longer branch/IBAN rows now wrap onto two lines, so give them extra height.
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30

# Leave the selection on the last filled merged block, as in the saved file.
$ws.Range("F11:F12").Select()
